# Append a new block of notes at the end of the document, right after the
# "cannot cancel requests" paragraph and before the pre-existing trailing
# blank paragraph.

$d = $word.ActiveDocument

# Index of the pre-existing trailing empty paragraph (just before sectPr).
$tailIndex = $d.Paragraphs.Count
$tailPara = $d.Paragraphs.Item($tailIndex)

# Insert all of the new paragraphs' text in one shot, right before that
# trailing blank paragraph. Each "`n" starts a new paragraph. The blank
# separator lines get a throwaway "Z" placeholder character so that their
# paragraph range isn't *just* the paragraph mark (which would make it
# impossible to later clear them out without merging with their neighbor).
$newText = "Wall app architecture`nActivity`nGoo.gl/z40GQJ`nZ`nSet maximum input length on textbox widget and save yourself a headache`nZ`n"
$tailPara.Range.InsertBefore($newText)

# The newly-inserted paragraphs now occupy indices starting at $tailIndex.
$wallPara     = $d.Paragraphs.Item($tailIndex)
$activityPara = $d.Paragraphs.Item($tailIndex + 1)
$linkPara     = $d.Paragraphs.Item($tailIndex + 2)
$blank1Para   = $d.Paragraphs.Item($tailIndex + 3)
$notePara     = $d.Paragraphs.Item($tailIndex + 4)
$blank2Para   = $d.Paragraphs.Item($tailIndex + 5)

# "Activity" becomes a single-level numbered list item using the
# "List Paragraph" style, like the other lists already in this document.
$activityPara.Style = "Paragrafoelenco"
$activityPara.Range.ListFormat.ApplyNumberDefault()
$lvl = $activityPara.Range.ListFormat.ListTemplate.ListLevels.Item(1)
$lvl.NumberFormat = "%1)"

# Strip the "Z" placeholders back out so the two separator paragraphs
# collapse to genuinely empty paragraphs, matching the rest of the doc.
# Go back-to-front so earlier character offsets stay valid.
$r2 = $d.Range($blank2Para.Range.Start, $blank2Para.Range.Start + 1)
$r2.Delete()
$r1 = $d.Range($blank1Para.Range.Start, $blank1Para.Range.Start + 1)
$r1.Delete()
